$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue $ws "D2" "69.464.75"
Set-TextValue $ws "E2" "  +0.03%  "
Set-TextValue $ws "D3" "2.492.64"
Set-TextValue $ws "E3" "  -0.80%  "
Set-TextValue $ws "E4" "  +0.04%  "
Set-TextValue $ws "D5" "568.77"
Set-TextValue $ws "E5" "  -0.63%  "
Set-TextValue $ws "D6" "166.05"
Set-TextValue $ws "E6" "  +0.07%  "
Set-TextValue $ws "E7" "  -0.01%  "
Set-TextValue $ws "E8" "  -0.31%  "
Set-TextValue $ws "E9" "  +0.60%  "
Set-TextValue $ws "E10" "  -0.74%  "
Set-TextValue $ws "D11" "0.347"
Set-TextValue $ws "E11" "  -2.94%  "
Set-TextValue $ws "D13" "2.949.17"
Set-TextValue $ws "E13" "  -0.81%  "
Set-TextValue $ws "D14" "69.363.70"
Set-TextValue $ws "E14" "  +0.04%  "
Set-TextValue $ws "E15" "  -0.27%  "
Set-TextValue $ws "D16" "24.15"
Set-TextValue $ws "E16" "  -2.62%  "
Set-TextValue $ws "D17" "2.483.41"
Set-TextValue $ws "E17" "  -1.85%  "
Set-TextValue $ws "D18" "11.19"
Set-TextValue $ws "E18" "  -0.74%  "
Set-TextValue $ws "E19" "  -1.44%  "
Set-TextValue $ws "D20" "352.94"
Set-TextValue $ws "E20" "  +1.19%  "
Set-TextValue $ws "E21" "  -0.02%  "
Set-TextValue $ws "D22" "1.91"
Set-TextValue $ws "E22" "  -3.32%  "
Set-TextValue $ws "E23" "  -0.08%  "
Set-TextValue $ws "D24" "69.35"
Set-TextValue $ws "E24" "  -1.31%  "
Set-TextValue $ws "D25" "3.79"
Set-TextValue $ws "E25" "  -2.84%  "
Set-TextValue $ws "E26" "  -0.87%  "
Set-TextValue $ws "E27" "  -2.46%  "
Set-TextValue $ws "E28" "  +0.15%  "
Set-TextValue $ws "E29" "  -1.63%  "
Set-TextValue $ws "D30" "7.53"
Set-TextValue $ws "E30" "  -3.83%  "
Set-TextValue $ws "D31" "3.50"
Set-TextValue $ws "E31" "  +136.06%  "
Set-TextValue $ws "E32" "  -3.28%  "
Set-TextValue $ws "D33" "438.77"
Set-TextValue $ws "E33" "  -4.83%  "
Set-TextValue $ws "D34" "1.00"
Set-TextValue $ws "E34" "  +0.06%  "
Set-TextValue $ws "E35" "  -0.82%  "
Set-TextValue $ws "E36" "  -3.40%  "
Set-TextValue $ws "D37" "153.08"
Set-TextValue $ws "E37" "  -2.57%  "
Set-TextValue $ws "D38" "19.06"
Set-TextValue $ws "E38" "  -0.08%  "
Set-TextValue $ws "D39" "18.13"
Set-TextValue $ws "E39" "  -1.69%  "
Set-TextValue $ws "E40" "  +0.01%  "
Set-TextValue $ws "D41" "0.313"
Set-TextValue $ws "E41" "  -1.21%  "
Set-TextValue $ws "D42" "4.58"
Set-TextValue $ws "E42" "  -2.16%  "
Set-TextValue $ws "E43" "  -1.90%  "
Set-TextValue $ws "E44" "  -2.43%  "
Set-TextValue $ws "E45" "  -3.70%  "
Set-TextValue $ws "D46" "139.15"
Set-TextValue $ws "E46" "  -1.82%  "
Set-TextValue $ws "E47" "  -0.84%  "
Set-TextValue $ws "E48" "  -2.63%  "
Set-TextValue $ws "E49" "  -1.00%  "
Set-TextValue $ws "E50" "  -0.71%  "
